$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BCI simultanea - update fecha_pago (B2) and nemotecnico (H2)
$ws.Range("B2").Value = 45621
$ws.Range("H2").Value = "SMT_20012025_25112024_0.52_BCI"

# Row 3: ORO BLANCO simultanea - update fecha_pago (B3) and nemotecnico (H3)
$ws.Range("B3").Value = 45621
$ws.Range("H3").Value = "SMT_23122024_25112024_0.56_ORO_BLANCO"

# Row 4: ORO BLANCO simultanea - update fecha_pago (B4) and nemotecnico (H4)
$ws.Range("B4").Value = 45621
$ws.Range("H4").Value = "SMT_23122024_25112024_0.56_ORO_BLANCO"
